# Re-applies a refreshed cryptos.com scrape onto the existing sheet:
# new Price (D) / Volume(1h) (E) figures for every coin row, plus the
# Uniswap / BitcoinCash rows trading rank places (row 22 <-> row 23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> refreshed value scraped this run.
$updates = [ordered]@{
    'D2' = '96.506.51'
    'E2' = '  -0.79%  '
    'D3' = '3.335.96'
    'E3' = '  -1.53%  '
    'E4' = '  -0.36%  '
    'D5' = '250.17'
    'E5' = '  -1.24%  '
    'D6' = '654.02'
    'E6' = '  +1.61%  '
    'E7' = '  -2.29%  '
    'D8' = '0.421'
    'E8' = '  +0.99%  '
    'E9' = '  -0.06%  '
    'D10' = '0.997'
    'E10' = '  -3.64%  '
    'D11' = '3.332.89'
    'E11' = '  -1.46%  '
    'D12' = '0.206'
    'E12' = '  -1.96%  '
    'D13' = '40.34'
    'E13' = '  -1.27%  '
    'D14' = '96.193.41'
    'E14' = '  -0.87%  '
    'D15' = '6.08'
    'E15' = '  -1.33%  '
    'D16' = '0.0000251'
    'E16' = '  -0.98%  '
    'D17' = '3.960.02'
    'E17' = '  -1.49%  '
    'D18' = '8.62'
    'E18' = '  +3.71%  '
    'D19' = '3.347.68'
    'E19' = '  -1.24%  '
    'D20' = '0.568'
    'E20' = '  +20.71%  '
    'D21' = '17.13'
    'E21' = '  +0.21%  '
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '10.59'
    'E22' = '  -0.63%  '
    'B23' = 'BitcoinCash'
    'C23' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D23' = '505.09'
    'E23' = '  +1.82%  '
    'D24' = '3.35'
    'E24' = '  -0.34%  '
    'E25' = '  -1.34%  '
    'D26' = '6.56'
    'E26' = '  +9.82%  '
    'D27' = '96.23'
    'E27' = '  +1.36%  '
    'D28' = '12.05'
    'E28' = '  -2.92%  '
    'E29' = '  -2.97%  '
    'E30' = '  +0.38%  '
    'D31' = '11.09'
    'E31' = '  -0.13%  '
    'E32' = '  -4.22%  '
    'D33' = '2.48'
    'E33' = '  +12.32%  '
    'E34' = '  +0.39%  '
    'D35' = '0.548'
    'E35' = '  -1.70%  '
    'D36' = '28.08'
    'E36' = '  -3.81%  '
    'D37' = '1.49'
    'E37' = '  +7.99%  '
    'D38' = '7.66'
    'E38' = '  +1.33%  '
    'E39' = '  +0.01%  '
    'E40' = '  +0.09%  '
    'D41' = '508.87'
    'E41' = '  +1.49%  '
    'D42' = '24.33'
    'E42' = '  -1.41%  '
    'D43' = '0.0429'
    'E43' = '  +6.02%  '
    'E44' = '  -1.55%  '
    'E45' = '  +0.44%  '
    'E46' = '  +8.76%  '
    'D47' = '5.53'
    'E47' = '  +2.70%  '
    'D48' = '8.41'
    'E48' = '  +4.57%  '
    'D49' = '53.44'
    'E49' = '  +3.99%  '
    'D50' = '3.12'
    'E50' = '  -1.42%  '
    'D51' = '162.85'
    'E51' = '  +1.64%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)
    # Plain-number-looking text (e.g. "96.23", "0.0000251") must stay text, just
    # like the original scrape output -- force Text format before writing so Excel
    # doesn't silently coerce it into a floating point number.
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.NumberFormat = '@'
    }
    $range.Value = $value
}

